{"js": "// Update the answer key table: each data row in the table holds 5\n// \"<dividend>\u00f7<divisor>=<quotient>, <remainder>\" answers; the rest of the\n// rows are blank spacer rows. We target cells by (row, col) rather than a\n// document-wide text search/replace because a couple of the old values\n// (\"38\u00f77=5, 3\") are duplicated but map to two different new values.\nconst edits = [\n  { row: 0, col: 0, oldText: \"78\u00f73=26, 0\", newText: \"40\u00f73=13, 1\" },\n  { row: 0, col: 1, oldText: \"88\u00f75=17, 3\", newText: \"80\u00f79=8, 8\" },\n  { row: 0, col: 2, oldText: \"77\u00f79=8, 5\", newText: \"30\u00f72=15, 0\" },\n  { row: 0, col: 3, oldText: \"23\u00f72=11, 1\", newText: \"66\u00f73=22, 0\" },\n  { row: 0, col: 4, oldText: \"12\u00f73=4, 0\", newText: \"69\u00f77=9, 6\" },\n  { row: 4, col: 0, oldText: \"68\u00f73=22, 2\", newText: \"23\u00f75=4, 3\" },\n  { row: 4, col: 1, oldText: \"62\u00f79=6, 8\", newText: \"73\u00f74=18, 1\" },\n  { row: 4, col: 2, oldText: \"61\u00f74=15, 1\", newText: \"90\u00f76=15, 0\" },\n  { row: 4, col: 3, oldText: \"21\u00f77=3, 0\", newText: \"61\u00f78=7, 5\" },\n  { row: 4, col: 4, oldText: \"19\u00f76=3, 1\", newText: \"80\u00f74=20, 0\" },\n  { row: 8, col: 0, oldText: \"75\u00f78=9, 3\", newText: \"35\u00f79=3, 8\" },\n  { row: 8, col: 1, oldText: \"60\u00f78=7, 4\", newText: \"96\u00f78=12, 0\" },\n  { row: 8, col: 2, oldText: \"50\u00f79=5, 5\", newText: \"36\u00f73=12, 0\" },\n  { row: 8, col: 3, oldText: \"69\u00f72=34, 1\", newText: \"24\u00f73=8, 0\" },\n  { row: 8, col: 4, oldText: \"38\u00f77=5, 3\", newText: \"70\u00f74=17, 2\" },\n  { row: 12, col: 0, oldText: \"67\u00f73=22, 1\", newText: \"42\u00f75=8, 2\" },\n  { row: 12, col: 1, oldText: \"38\u00f77=5, 3\", newText: \"46\u00f78=5, 6\" },\n  { row: 12, col: 2, oldText: \"10\u00f73=3, 1\", newText: \"49\u00f73=16, 1\" },\n  { row: 12, col: 3, oldText: \"73\u00f73=24, 1\", newText: \"38\u00f76=6, 2\" },\n  { row: 12, col: 4, oldText: \"16\u00f79=1, 7\", newText: \"29\u00f73=9, 2\" },\n  { row: 16, col: 0, oldText: \"68\u00f76=11, 2\", newText: \"58\u00f76=9, 4\" },\n  { row: 16, col: 1, oldText: \"66\u00f77=9, 3\", newText: \"97\u00f73=32, 1\" },\n  { row: 16, col: 2, oldText: \"27\u00f76=4, 3\", newText: \"49\u00f73=16, 1\" },\n  { row: 16, col: 3, oldText: \"34\u00f78=4, 2\", newText: \"96\u00f76=16, 0\" },\n  { row: 16, col: 4, oldText: \"51\u00f78=6, 3\", newText: \"52\u00f76=8, 4\" },\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\nconst table = tables.items[0];\n\n// Load every cell's current text first so we can confirm we are editing\n// the expected cell before mutating it.\nconst cells = edits.map((e) => table.getCell(e.row, e.col));\ncells.forEach((c) => c.load(\"value\"));\nawait context.sync();\n\nfor (let i = 0; i < edits.length; i++) {\n  const e = edits[i];\n  const cell = cells[i];\n  const current = cell.value;\n  if (current !== e.oldText) {\n    throw new Error(\n      `Unexpected text at row ${e.row}, col ${e.col}: \"${current}\" (expected \"${e.oldText}\")`\n    );\n  }\n  // Replace just the text of the cell's range so the existing run\n  // formatting (font, size, paragraph alignment) is preserved.\n  const range = cell.body.getRange();\n  range.insertText(e.newText, Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# Update the answer key table: each data row in the table holds 5\n# \"<dividend>\u00f7<divisor>=<quotient>, <remainder>\" answers; the rest of the\n# rows are blank spacer rows. We target cells by (row, col) rather than a\n# document-wide Find/Replace because a couple of the old values\n# (\"38\u00f77=5, 3\") are duplicated but map to two different new values.\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$edits = @(\n  @{ Row = 1; Col = 1; OldText = \"78\u00f73=26, 0\"; NewText = \"40\u00f73=13, 1\" },\n  @{ Row = 1; Col = 2; OldText = \"88\u00f75=17, 3\"; NewText = \"80\u00f79=8, 8\" },\n  @{ Row = 1; Col = 3; OldText = \"77\u00f79=8, 5\"; NewText = \"30\u00f72=15, 0\" },\n  @{ Row = 1; Col = 4; OldText = \"23\u00f72=11, 1\"; NewText = \"66\u00f73=22, 0\" },\n  @{ Row = 1; Col = 5; OldText = \"12\u00f73=4, 0\"; NewText = \"69\u00f77=9, 6\" },\n  @{ Row = 5; Col = 1; OldText = \"68\u00f73=22, 2\"; NewText = \"23\u00f75=4, 3\" },\n  @{ Row = 5; Col = 2; OldText = \"62\u00f79=6, 8\"; NewText = \"73\u00f74=18, 1\" },\n  @{ Row = 5; Col = 3; OldText = \"61\u00f74=15, 1\"; NewText = \"90\u00f76=15, 0\" },\n  @{ Row = 5; Col = 4; OldText = \"21\u00f77=3, 0\"; NewText = \"61\u00f78=7, 5\" },\n  @{ Row = 5; Col = 5; OldText = \"19\u00f76=3, 1\"; NewText = \"80\u00f74=20, 0\" },\n  @{ Row = 9; Col = 1; OldText = \"75\u00f78=9, 3\"; NewText = \"35\u00f79=3, 8\" },\n  @{ Row = 9; Col = 2; OldText = \"60\u00f78=7, 4\"; NewText = \"96\u00f78=12, 0\" },\n  @{ Row = 9; Col = 3; OldText = \"50\u00f79=5, 5\"; NewText = \"36\u00f73=12, 0\" },\n  @{ Row = 9; Col = 4; OldText = \"69\u00f72=34, 1\"; NewText = \"24\u00f73=8, 0\" },\n  @{ Row = 9; Col = 5; OldText = \"38\u00f77=5, 3\"; NewText = \"70\u00f74=17, 2\" },\n  @{ Row = 13; Col = 1; OldText = \"67\u00f73=22, 1\"; NewText = \"42\u00f75=8, 2\" },\n  @{ Row = 13; Col = 2; OldText = \"38\u00f77=5, 3\"; NewText = \"46\u00f78=5, 6\" },\n  @{ Row = 13; Col = 3; OldText = \"10\u00f73=3, 1\"; NewText = \"49\u00f73=16, 1\" },\n  @{ Row = 13; Col = 4; OldText = \"73\u00f73=24, 1\"; NewText = \"38\u00f76=6, 2\" },\n  @{ Row = 13; Col = 5; OldText = \"16\u00f79=1, 7\"; NewText = \"29\u00f73=9, 2\" },\n  @{ Row = 17; Col = 1; OldText = \"68\u00f76=11, 2\"; NewText = \"58\u00f76=9, 4\" },\n  @{ Row = 17; Col = 2; OldText = \"66\u00f77=9, 3\"; NewText = \"97\u00f73=32, 1\" },\n  @{ Row = 17; Col = 3; OldText = \"27\u00f76=4, 3\"; NewText = \"49\u00f73=16, 1\" },\n  @{ Row = 17; Col = 4; OldText = \"34\u00f78=4, 2\"; NewText = \"96\u00f76=16, 0\" },\n  @{ Row = 17; Col = 5; OldText = \"51\u00f78=6, 3\"; NewText = \"52\u00f76=8, 4\" }\n)\n\nforeach ($edit in $edits) {\n    $cell = $t.Cell($edit.Row, $edit.Col)\n    $range = $cell.Range\n    $current = $range.Text.TrimEnd([char]13, [char]7)\n    if ($current -ne $edit.OldText) {\n        throw \"Unexpected text at row $($edit.Row), col $($edit.Col): '$current' (expected '$($edit.OldText)')\"\n    }\n    # Assigning Range.Text replaces just the cell's text run, preserving\n    # the existing run formatting (font, size, paragraph alignment).\n    $range.Text = $edit.NewText\n}\n"}
